$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F9").Value = 25059
$wsExhibit.Range("F11").Value = 232
$wsExhibit.Range("F16").Value = 333
$wsExhibit.Range("F17").Value = 185
$wsExhibit.Range("F18").Value = 159
$wsExhibit.Range("F20").Value = 253
$wsExhibit.Range("F23").Value = 1412

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F10").Value = 406
$wsShow.Range("F15").Value = 18
$wsShow.Range("F17").Value = 22

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 4765

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 4765
$wsAll.Range("F12").Value = 25059
$wsAll.Range("F15").Value = 232
$wsAll.Range("F27").Value = 406
$wsAll.Range("F31").Value = 333
$wsAll.Range("F32").Value = 185
$wsAll.Range("F33").Value = 159
$wsAll.Range("F36").Value = 253
$wsAll.Range("F41").Value = 18
$wsAll.Range("F42").Value = 1412
$wsAll.Range("F48").Value = 22
